$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 2155
$ws.Range("B2").Value = 2589
$ws.Range("C2").Value = 4283
$ws.Range("D2").Value = 3686
$ws.Range("E2").Value = 3483
$ws.Range("F2").Value = 2168
